# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Sat May 13 05:05:10 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.923.71"
$ws.Range("E2").Value = "  +0.43%  "

# Row 3
$ws.Range("D3").Value = "1.810.86"
$ws.Range("E3").Value = "  +1.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'312.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

# Row 6
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("E7").Value = "  -2.75%  "

# Row 8
$ws.Range("D8").Value = "'0.3680"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.30%  "

# Row 9
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("D10").Value = "'0.8616"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.70%  "

# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'21.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.29%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "2.025.86"
$ws.Range("E12").Value = "  +12.35%  "

# Row 13
$ws.Range("D13").Value = "'6.619"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.05%  "

# Row 14
$ws.Range("D14").Value = "'5.383"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.14%  "

# Row 15
$ws.Range("D15").Value = "'0.06903"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16
$ws.Range("D16").Value = "'80.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.32%  "

# Row 17
$ws.Range("D17").Value = "'1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("D18").Value = "'0.000008915"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.40%  "

# Row 19
$ws.Range("D19").Value = "'1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "

# Row 20
$ws.Range("D20").Value = "'15.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "

# Row 21
$ws.Range("D21").Value = "26.964.85"
$ws.Range("E21").Value = "  +0.25%  "

# Row 22
$ws.Range("D22").Value = "'5.180"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "

# Row 23
$ws.Range("D23").Value = "'10.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.07%  "

# Row 24
$ws.Range("D24").Value = "2.239.29"
$ws.Range("E24").Value = "  +10.97%  "

# Row 25
$ws.Range("D25").Value = "'153.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "

# Row 26
$ws.Range("D26").Value = "'1.884"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.75%  "

# Row 27
$ws.Range("D27").Value = "'18.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "

# Row 28
$ws.Range("E28").Value = "  +3.14%  "

# Row 29
$ws.Range("D29").Value = "'1.888"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.61%  "

# Row 30
$ws.Range("D30").Value = "'114.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "

# Row 31
$ws.Range("D31").Value = "'0.08925"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.12%  "

# Row 32
$ws.Range("D32").Value = "'0.7405"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.27%  "

# Row 33
$ws.Range("E33").Value = "  +5.50%  "

# Row 34
$ws.Range("D34").Value = "'4.415"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.83%  "

# Row 35
$ws.Range("D35").Value = "'2.802"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "

# Row 36
$ws.Range("D36").Value = "'1.009"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "

# Row 37
$ws.Range("D37").Value = "'1.123"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.45%  "

# Row 38
$ws.Range("D38").Value = "'0.05208"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.24%  "

# Row 39
$ws.Range("D39").Value = "'0.01918"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.32%  "

# Row 40
$ws.Range("D40").Value = "'0.5078"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.12%  "

# Row 41
$ws.Range("D41").Value = "'2.750"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.55%  "

# Row 42
$ws.Range("D42").Value = "'0.1643"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.96%  "

# Row 43
$ws.Range("D43").Value = "'6.400"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.75%  "

# Row 44
$ws.Range("D44").Value = "'8.228"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.60%  "

# Row 45
$ws.Range("D45").Value = "'106.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.83%  "

# Row 46
$ws.Range("D46").Value = "'10.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.13%  "

# Row 47
$ws.Range("E47").Value = "  -0.09%  "

# Row 48
$ws.Range("E48").Value = "  +4.49%  "

# Row 49
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.4553"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.45%  "

# Row 51
$ws.Range("D51").Value = "'1.796"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.81%  "
